# Replace the 4-bullet "KEY ACHIEVEMENTS AND IMPACT" list with the new
# 6-bullet list that includes the Supreme Court mentions, per the commit:
#   "Fix Supreme Court mention missing from short data engineering resumes"
#
# Before (4 bullets):
#   • Discovered systematic race coding errors affecting all Black and Asian-American voters
#   • Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M
#   • Built redistricting platform used by thousands of analysts nationwide
#   • Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%
#
# After (6 bullets):
#   • Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs **73.5%**
#   • **$4.7M** savings enabled nonprofit access
#   • Legal precedent: Data analysis utilized in Supreme Court case
#   • Expert methodology validated at highest judicial level
#   • Breakthrough demographic discovery: Uncovered systematic voter miscoding affecting millions
#   • **178%** accuracy improvement in racial classification algorithms

$d = $word.ActiveDocument

$CR = [char]13
$BEL = [char]7
$boldColor = 5258796   # BGR-packed 0x2C3E50 ("2C3E50" swatch used throughout the doc)

function Clean-Text($t) {
    return $t.TrimEnd($CR, $BEL)
}

# --- Locate the short "Discovered ..." bullet that lives inside the
# --- "KEY ACHIEVEMENTS AND IMPACT" section (the Siege Analytics experience
# --- bullet has the same opening words but extra trailing text, so match
# --- on exact equality to disambiguate; if there are duplicates use the
# --- last one, which is the short/achievements-section copy).
$anchorText = "• Discovered systematic race coding errors affecting all Black and Asian-American voters"
$algoText   = "• Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M"
$builtText  = "• Built redistricting platform used by thousands of analysts nationwide"
$achText    = "• Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%"

$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = Clean-Text $d.Paragraphs($i).Range.Text
    if ($t -eq $anchorText) {
        $candidate = $i
        $n2 = Clean-Text $d.Paragraphs($i + 1).Range.Text
        $n3 = Clean-Text $d.Paragraphs($i + 2).Range.Text
        $n4 = Clean-Text $d.Paragraphs($i + 3).Range.Text
        if ($n2 -eq $algoText -and $n3 -eq $builtText -and $n4 -eq $achText) {
            $anchorIndex = $candidate
        }
    }
}
if ($anchorIndex -eq -1) {
    throw "Could not locate the KEY ACHIEVEMENTS AND IMPACT bullet block"
}

# --- Insert 6 empty paragraphs right after the (still plain, non-bold)
# --- anchor paragraph. Doing all the InsertParagraphAfter calls off the
# --- same still-unedited anchor range (rather than chaining off each new
# --- paragraph) keeps every new paragraph's inherited run formatting
# --- clean/non-bold.
$anchorRange = $d.Paragraphs($anchorIndex).Range
$anchorRange.InsertParagraphAfter()
$anchorRange.InsertParagraphAfter()
$anchorRange.InsertParagraphAfter()
$anchorRange.InsertParagraphAfter()
$anchorRange.InsertParagraphAfter()
$anchorRange.InsertParagraphAfter()

# --- Bullet 1: "Algorithmic innovation ... **73.5%**" (bold 73.5%)
$p1 = $d.Paragraphs($anchorIndex + 1)
$p1.Range.Text = "• Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs **73.5%**"
$full = $p1.Range.Text
$idx = $full.IndexOf("73.5%")
$s = $p1.Range.Start + $idx
$e = $s + 5
$rng = $d.Range($s, $e)
$rng.Font.Bold = 1
$rng.Font.Color = $boldColor

# --- Bullet 2: "**$4.7M** savings enabled nonprofit access" (bold $4.7M)
$p2 = $d.Paragraphs($anchorIndex + 2)
$p2.Range.Text = "• **`$4.7M** savings enabled nonprofit access"
$full = $p2.Range.Text
$idx = $full.IndexOf("`$4.7M")
$s = $p2.Range.Start + $idx
$e = $s + 5
$rng = $d.Range($s, $e)
$rng.Font.Bold = 1
$rng.Font.Color = $boldColor

# --- Bullet 3: Legal precedent (plain text)
$p3 = $d.Paragraphs($anchorIndex + 3)
$p3.Range.Text = "• Legal precedent: Data analysis utilized in Supreme Court case"

# --- Bullet 4: Expert methodology (plain text)
$p4 = $d.Paragraphs($anchorIndex + 4)
$p4.Range.Text = "• Expert methodology validated at highest judicial level"

# --- Bullet 5: Breakthrough demographic discovery (plain text)
$p5 = $d.Paragraphs($anchorIndex + 5)
$p5.Range.Text = "• Breakthrough demographic discovery: Uncovered systematic voter miscoding affecting millions"

# --- Bullet 6: "**178%** accuracy improvement ..." (bold 178%)
$p6 = $d.Paragraphs($anchorIndex + 6)
$p6.Range.Text = "• **178%** accuracy improvement in racial classification algorithms"
$full = $p6.Range.Text
$idx = $full.IndexOf("178%")
$s = $p6.Range.Start + $idx
$e = $s + 4
$rng = $d.Range($s, $e)
$rng.Font.Bold = 1
$rng.Font.Color = $boldColor

# --- Now remove the 4 original bullets. The first (anchor) paragraph is
# --- still sitting right before our 6 new ones; the other 3 original
# --- bullets got pushed down to directly follow the 6 new paragraphs.
# --- Delete the trailing trio first (doesn't disturb the anchor index),
# --- then delete the anchor paragraph itself.
$oldTailStart = $d.Paragraphs($anchorIndex + 7)
$oldTailEnd   = $d.Paragraphs($anchorIndex + 9)
$tailRange = $d.Range($oldTailStart.Range.Start, $oldTailEnd.Range.End)
$tailRange.Delete()

$anchorParaRange = $d.Paragraphs($anchorIndex).Range
$anchorParaRange.Delete()

Write-Output "Done. Final paragraph count: $($d.Paragraphs.Count)"
